# Apply odds-update edits to the weekly FlashScore games sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 ---
$ws.Range("G8").Value = 3.85
$ws.Range("H8").Value = 3.55
$ws.Range("P8").Value = 1.37
$ws.Range("R8").Value = 1.6
$ws.Range("T8").Value = 12.5
$ws.Range("U8").Value = 26
$ws.Range("AA8").Value = 7.3
$ws.Range("AI8").Value = 23

# --- Row 25 ---
$ws.Range("G25").Value = 1.25
$ws.Range("I25").Value = 9

# --- Row 28 ---
$ws.Range("H28").Value = 3.2
$ws.Range("I28").Value = 3.05
$ws.Range("R28").Value = 1.9
$ws.Range("S28").Value = 1.72
$ws.Range("AA28").Value = 6.2
$ws.Range("AD28").Value = 8

# --- Row 37 ---
$ws.Range("G37").Value = 2.37
$ws.Range("H37").Value = 2.92
$ws.Range("I37").Value = 3.05
$ws.Range("L37").Value = 1.44
$ws.Range("M37").Value = 2.42
$ws.Range("N37").Value = 2.25
$ws.Range("O37").Value = 1.5
$ws.Range("P37").Value = 1.5
$ws.Range("Q37").Value = 2.27
$ws.Range("R37").Value = 1.91
$ws.Range("S37").Value = 1.7
$ws.Range("T37").Value = 6.2
$ws.Range("U37").Value = 10.5
$ws.Range("V37").Value = 9.75
$ws.Range("W37").Value = 25
$ws.Range("X37").Value = 24
$ws.Range("Y37").Value = 40
$ws.Range("Z37").Value = 6.9
$ws.Range("AB37").Value = 16
$ws.Range("AC37").Value = 90
$ws.Range("AD37").Value = 7.8
$ws.Range("AF37").Value = 11.25
$ws.Range("AH37").Value = 30
$ws.Range("AI37").Value = 45
$ws.Range("AJ37").Value = 900
